# Automatic update of files.
#
# 1) Every "Förändrad" (column C) timestamp on a data row is refreshed to
#    the new export date (45184 -> 45186).
# 2) Every HYPERLINK(...) formula (artfynd/kartor/knärot/klagomål/
#    klagomålsmail/tillsyn/tillsynsmail columns) gains a second argument:
#    the friendly display text, which is the row's "Beteckning" (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 387
$firstDataRow = 2
$firstCol = 1
$lastCol = 25
$newDate = 45186

for ($row = $firstDataRow; $row -le $lastRow; $row++) {

    $beteckning = $ws.Cells.Item($row, 1).Value()

    # --- update the "Förändrad" date column (C = column 3) ---
    $cCell = $ws.Cells.Item($row, 3)
    if ($cCell.Value() -ne $null) {
        $cCell.Value = $newDate
    }

    # --- add the friendly-name argument to every HYPERLINK formula on the row ---
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        if ($cell.HasFormula()) {
            $formula = $cell.Formula()
            if ($formula.IndexOf("HYPERLINK(") -ge 0 -and $formula.IndexOf(",") -lt 0) {
                $trimmed = $formula.TrimEnd()
                $withoutParen = $trimmed.Substring(0, $trimmed.Length - 1)
                $escaped = $beteckning.Replace('"', '""')
                $newFormula = $withoutParen + ', "' + $escaped + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
